$d = $word.ActiveDocument

# --- Part 1: "Musica" paragraph -----------------------------------------
# The first paragraph currently holds a gramStart/gramEnd proofErr pair
# around a "Musica" run, followed by a separate " " run:
#   <w:proofErr w:type="gramStart"/><w:r><w:t>Musica</w:t></w:r>
#   <w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r>
# Target: a single run "Musica " with no proofErr markers.
# Deleting the whole paragraph range (content + paragraph mark) drops the
# proofErr markers along with it; re-inserting a fresh paragraph in front
# of the (now first) paragraph and typing the text gives a clean single run.
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Delete()

$newFirstAnchor = $d.Paragraphs.Item(1).Range
$newFirstAnchor.InsertParagraphBefore()
$d.Paragraphs.Item(1).Range.Text = "Musica "

# --- Part 2: add a new "Xuxa" paragraph ----------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$d.Paragraphs.Item($d.Paragraphs.Count).Range.Text = "Xuxa"
